# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4370
$ws1.Range("F6").Value = 50
$ws1.Range("C7").Value = "南宁·排球少年ONLY（取消）"
$ws1.Range("G7").Value = "不可售"
$ws1.Range("F8").Value = 216
$ws1.Range("F9").Value = 131
$ws1.Range("F11").Value = 161
$ws1.Range("F12").Value = 1630
$ws1.Range("F13").Value = 297
$ws1.Range("F14").Value = 3451
$ws1.Range("F15").Value = 230

# ---- Sheet: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 43

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4370
$ws4.Range("F7").Value = 50
$ws4.Range("C8").Value = "南宁·排球少年ONLY（取消）"
$ws4.Range("G8").Value = "不可售"
$ws4.Range("F9").Value = 43
$ws4.Range("F10").Value = 216
$ws4.Range("F11").Value = 131
$ws4.Range("F13").Value = 161
$ws4.Range("F16").Value = 1630
$ws4.Range("F17").Value = 297
$ws4.Range("F18").Value = 3451
$ws4.Range("F19").Value = 230
